# Update the "Förändrad" (Changed) date column C for rows 2-10 from
# 45184 (2023-09-15) to 45185 (2023-09-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45185
    }
}
